$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format so numeric-looking strings
# such as "1.000" or "321.35" are preserved verbatim instead of being
# coerced into numbers by Excel's input parser.
$ws.Range("D2:D51").NumberFormat = "@"


# Row 2
$ws.Range("D2").Value = '28.712.82'
$ws.Range("E2").Value = '  -2.66%  '

# Row 3
$ws.Range("D3").Value = '1.883.41'
$ws.Range("E3").Value = '  -5.25%  '

# Row 4
$ws.Range("E4").Value = '  +0.35%  '

# Row 5
$ws.Range("D5").Value = '321.35'
$ws.Range("E5").Value = '  -1.56%  '

# Row 6
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  +0.38%  '

# Row 7
$ws.Range("D7").Value = '0.4555'
$ws.Range("E7").Value = '  -1.82%  '

# Row 8
$ws.Range("D8").Value = '0.3786'
$ws.Range("E8").Value = '  -4.24%  '

# Row 9
$ws.Range("D9").Value = '45.44'
$ws.Range("E9").Value = '  -1.72%  '

# Row 10
$ws.Range("D10").Value = '0.07697'
$ws.Range("E10").Value = '  -2.75%  '

# Row 11
$ws.Range("D11").Value = '0.9587'
$ws.Range("E11").Value = '  -4.45%  '

# Row 12
$ws.Range("D12").Value = '21.93'
$ws.Range("E12").Value = '  -2.46%  '

# Row 13
$ws.Range("D13").Value = '1.877.83'
$ws.Range("E13").Value = '  -4.96%  '

# Row 14
$ws.Range("D14").Value = '6.919'
$ws.Range("E14").Value = '  -4.31%  '

# Row 15
$ws.Range("D15").Value = '5.634'
$ws.Range("E15").Value = '  -3.91%  '

# Row 16
$ws.Range("D16").Value = '0.06986'
$ws.Range("E16").Value = '  -1.89%  '

# Row 17
$ws.Range("E17").Value = '  +0.43%  '

# Row 18
$ws.Range("D18").Value = '82.73'
$ws.Range("E18").Value = '  -7.07%  '

# Row 19
$ws.Range("D19").Value = '0.000009473'
$ws.Range("E19").Value = '  -5.14%  '

# Row 20
$ws.Range("D20").Value = '16.55'
$ws.Range("E20").Value = '  -3.78%  '

# Row 21
$ws.Range("D21").Value = '1.000'
$ws.Range("E21").Value = '  +0.41%  '

# Row 22
$ws.Range("D22").Value = '28.695.55'
$ws.Range("E22").Value = '  -3.01%  '

# Row 23
$ws.Range("D23").Value = '5.294'
$ws.Range("E23").Value = '  -5.08%  '

# Row 24
$ws.Range("D24").Value = '10.83'
$ws.Range("E24").Value = '  -3.77%  '

# Row 25
$ws.Range("D25").Value = '2.111.68'
$ws.Range("E25").Value = '  -4.72%  '

# Row 26
$ws.Range("D26").Value = '2.072'
$ws.Range("E26").Value = '  -2.17%  '

# Row 27
$ws.Range("D27").Value = '154.53'
$ws.Range("E27").Value = '  -2.10%  '

# Row 28
$ws.Range("D28").Value = '18.89'
$ws.Range("E28").Value = '  -4.08%  '

# Row 29
$ws.Range("D29").Value = '5.591'
$ws.Range("E29").Value = '  -7.29%  '

# Row 30
$ws.Range("D30").Value = '116.44'
$ws.Range("E30").Value = '  -3.51%  '

# Row 31
$ws.Range("D31").Value = '1.800'
$ws.Range("E31").Value = '  -6.33%  '

# Row 32
$ws.Range("D32").Value = '0.09218'
$ws.Range("E32").Value = '  -2.05%  '

# Row 33
$ws.Range("D33").Value = '0.8401'
$ws.Range("E33").Value = '  -6.23%  '

# Row 34
$ws.Range("D34").Value = '5.042'
$ws.Range("E34").Value = '  -4.56%  '

# Row 35
$ws.Range("D35").Value = '1.237'
$ws.Range("E35").Value = '  -8.54%  '

# Row 36
$ws.Range("D36").Value = '2.979'
$ws.Range("E36").Value = '  -5.84%  '

# Row 37
$ws.Range("D37").Value = '0.05642'
$ws.Range("E37").Value = '  -3.31%  '

# Row 38
$ws.Range("D38").Value = '1.140'
$ws.Range("E38").Value = '  -3.31%  '

# Row 39
$ws.Range("D39").Value = '1.000'
$ws.Range("E39").Value = '  +0.42%  '

# Row 40
$ws.Range("D40").Value = '0.02020'
$ws.Range("E40").Value = '  -5.18%  '

# Row 41
$ws.Range("D41").Value = '7.419'
$ws.Range("E41").Value = '  -6.24%  '

# Row 42
$ws.Range("D42").Value = '0.5465'
$ws.Range("E42").Value = '  -5.44%  '

# Row 43
$ws.Range("D43").Value = '0.000003026'
$ws.Range("E43").Value = '  -23.76%  '

# Row 44
$ws.Range("D44").Value = '0.1740'
$ws.Range("E44").Value = '  -4.63%  '

# Row 45
$ws.Range("D45").Value = '9.126'
$ws.Range("E45").Value = '  -7.07%  '

# Row 46
$ws.Range("D46").Value = '2.670'
$ws.Range("E46").Value = '  +1.35%  '

# Row 47
$ws.Range("D47").Value = '0.5130'
$ws.Range("E47").Value = '  -4.62%  '

# Row 48
$ws.Range("D48").Value = '11.18'
$ws.Range("E48").Value = '  -8.35%  '

# Row 49
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '0.06779'
$ws.Range("E49").Value = '  -3.05%  '

# Row 50
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").Value = '2.053'
$ws.Range("E50").Value = '  -6.93%  '

# Row 51
$ws.Range("D51").Value = '116.45'
$ws.Range("E51").Value = '  -3.33%  '
